$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.72 = 14197.1 pesos`n✅ 14197.1 pesos = 3.69 = 921.89 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the "tasas" sheet numeric rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 269.1
$ws2.Range("O10").Value = 3820.44
$ws2.Range("N12").Value = 3850
$ws2.Range("O12").Value = 250.001
